$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("G2").Value = 108.9481836666667
$ws.Range("H2").Value = 326.844551
$ws.Range("I2").Value = 0.1523660837152667
$ws.Range("J2").Value = 0.1650457680857909
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09302566666666667
$ws.Range("N2").Value = 0.279077
$ws.Range("Q2").Value = 10.13497741771411
$ws.Range("R2").Value = 91.21479675942702
$ws.Range("S2").Value = 0.1523660837152667
$ws.Range("T2").Value = 0.1650457680857909

# Row 3
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("I3").Value = 0.2954065074566193
$ws.Range("J3").Value = 0.3199898083081954
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.09302566666666667
$ws.Range("N3").Value = 0.279077
$ws.Range("Q3").Value = 19.64963730191778
$ws.Range("R3").Value = 176.84673571726
$ws.Range("S3").Value = 0.2954065074566193
$ws.Range("T3").Value = 0.3199898083081954

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 109.1710686666667
$ws.Range("H4").Value = 327.513206
$ws.Range("I4").Value = 0.1526777925792968
$ws.Range("J4").Value = 0.1653834169091284
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09302566666666667
$ws.Range("N4").Value = 0.279077
$ws.Range("Q4").Value = 10.15571144342911
$ws.Range("R4").Value = 91.40140299086201
$ws.Range("S4").Value = 0.1526777925792968
$ws.Range("T4").Value = 0.1653834169091284

# Row 5
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("G5").Value = 164.799919
$ws.Range("H5").Value = 329.599838
$ws.Range("I5").Value = 0.2304757859153342
$ws.Range("J5").Value = 0.166437097565877
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.09302566666666667
$ws.Range("N5").Value = 0.279077
$ws.Range("Q5").Value = 15.33062233158767
$ws.Range("R5").Value = 91.983733989526
$ws.Range("S5").Value = 0.2304757859153342
$ws.Range("T5").Value = 0.166437097565877

# Row 6
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("G6").Value = 120.894928
$ws.Range("H6").Value = 362.684784
$ws.Range("I6").Value = 0.169073830333483
$ws.Range("J6").Value = 0.1831439091310082
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.09302566666666667
$ws.Range("N6").Value = 0.279077
$ws.Range("Q6").Value = 11.24633127381867
$ws.Range("R6").Value = 101.216981464368
$ws.Range("S6").Value = 0.169073830333483
$ws.Range("T6").Value = 0.1831439091310082
